# Split the "M=0.44 & catch" run on the "Season 2" slide into two runs:
#   "M=0.44 & " and "catch (note that this differs from 2015 model)"
# (minor edits to SMBKC document)

$p = $ppt.ActivePresentation

$slide = $null
$shape = $null

# Locate the slide/shape that contains the text we need to edit, rather than
# relying on a hard-coded slide index.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $shp = $sl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*M=0.44 & catch*") {
                $slide = $sl
                $shape = $shp
            }
        }
    }
}

$tr = $shape.TextFrame.TextRange

# Find the paragraph containing the target phrase and locate "catch" within it.
# (Paragraphs(...).Text carries a trailing carriage-return character, so use
# -like rather than an exact -eq match.)
$targetPara = $null
for ($k = 1; $k -le $tr.Paragraphs().Count; $k++) {
    $para = $tr.Paragraphs($k, 1)
    if ($para.Text -like "M=0.44 & catch*") {
        $targetPara = $para
    }
}

$found = $targetPara.Find("catch")

# Replace just the "catch" substring with the longer phrase; PowerPoint's
# text engine keeps the untouched leading text "M=0.44 & " as its own run
# and creates a new run for the replacement text, exactly matching the
# target split.
$sub = $tr.Characters($found.Start, $found.Length)
$sub.Text = "catch (note that this differs from 2015 model)"
